$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2612.4
$ws.Range("I98").Value = 1582.409
$ws.Range("J98").Value = 10165.667
$ws.Range("K98").Value = 1582.409
$ws.Range("L98").Value = 10165.667
$ws.Range("M98").Value = -84.40900000000011
$ws.Range("N98").Value = -13161.667

$ws.Range("H122").Value = 2612.4
$ws.Range("I122").Value = 1582.409
$ws.Range("J122").Value = 10165.667
$ws.Range("K122").Value = 4747.227000000001
$ws.Range("L122").Value = 30497.001
$ws.Range("M122").Value = -2297.227000000001
$ws.Range("N122").Value = -35397.001

$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 9000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -13920

$ws.Range("H137").Value = 1276.2858
$ws.Range("I137").Value = 1001
$ws.Range("J137").Value = 1386.4
$ws.Range("K137").Value = 3003
$ws.Range("L137").Value = 4159.200000000001
$ws.Range("M137").Value = -453
$ws.Range("N137").Value = -9259.200000000001

$ws.Range("H138").Value = 4998.8
$ws.Range("I138").Value = 3000
$ws.Range("J138").Value = 5498.5
$ws.Range("K138").Value = 9000
$ws.Range("L138").Value = 16495.5
$ws.Range("M138").Value = -3860
$ws.Range("N138").Value = -26775.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2079
$ws.Range("I45").Value = 2332.3333
$ws.Range("J45").Value = 1699
$ws.Range("K45").Value = 2332.3333
$ws.Range("L45").Value = 1699
$ws.Range("M45").Value = -1955.3333
$ws.Range("N45").Value = -2453

$ws.Range("H61").Value = 6799.294
$ws.Range("I61").Value = 5999.5
$ws.Range("J61").Value = 6905.933
$ws.Range("K61").Value = 5999.5
$ws.Range("L61").Value = 6905.933
$ws.Range("M61").Value = -5787.5
$ws.Range("N61").Value = -7329.933

$ws.Range("H63").Value = 2382.6562
$ws.Range("I63").Value = 2358.1667
$ws.Range("K63").Value = 2358.1667
$ws.Range("M63").Value = -1672.1667

$ws.Range("H66").Value = 2382.6562
$ws.Range("I66").Value = 2358.1667
$ws.Range("K66").Value = 11790.8335
$ws.Range("M66").Value = -8358.833500000001

$ws.Range("H74").Value = 2771.3333
$ws.Range("I74").Value = 2467.7
$ws.Range("J74").Value = 3378.6
$ws.Range("K74").Value = 2467.7
$ws.Range("L74").Value = 3378.6
$ws.Range("M74").Value = -1593.7
$ws.Range("N74").Value = -5126.6

$ws.Range("H77").Value = 2771.3333
$ws.Range("I77").Value = 2467.7
$ws.Range("J77").Value = 3378.6
$ws.Range("K77").Value = 12338.5
$ws.Range("L77").Value = 16893
$ws.Range("M77").Value = -7970.5
$ws.Range("N77").Value = -25629

$ws.Range("H136").Value = 6799.294
$ws.Range("I136").Value = 5999.5
$ws.Range("J136").Value = 6905.933
$ws.Range("K136").Value = 17998.5
$ws.Range("L136").Value = 20717.799
$ws.Range("M136").Value = -15448.5
$ws.Range("N136").Value = -25817.799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6264.7896
$ws.Range("I31").Value = 2319.45
$ws.Range("J31").Value = 8397.405000000001
$ws.Range("K31").Value = 2319.45
$ws.Range("L31").Value = 8397.405000000001
$ws.Range("M31").Value = -2024.45
$ws.Range("N31").Value = -8987.405000000001

$ws.Range("H34").Value = 6264.7896
$ws.Range("I34").Value = 2319.45
$ws.Range("J34").Value = 8397.405000000001
$ws.Range("K34").Value = 2319.45
$ws.Range("L34").Value = 8397.405000000001
$ws.Range("M34").Value = -2117.45
$ws.Range("N34").Value = -8801.405000000001

$ws.Range("H132").Value = 3373.4546
$ws.Range("I132").Value = 3026.125
$ws.Range("J132").Value = 4299.6665
$ws.Range("K132").Value = 9078.375
$ws.Range("L132").Value = 12898.9995
$ws.Range("M132").Value = -6548.375
$ws.Range("N132").Value = -17958.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 26316216
$ws.Range("I4").Value = 33333676
$ws.Range("J4").Value = 743.375
$ws.Range("K4").Value = 100001028
$ws.Range("L4").Value = 2230.125
$ws.Range("M4").Value = -100000916
$ws.Range("N4").Value = -2454.125

$ws.Range("H68").Value = 3389.0454
$ws.Range("J68").Value = 3431.3809
$ws.Range("L68").Value = 10294.1427
$ws.Range("N68").Value = -11916.1427

$ws.Range("H71").Value = 3389.0454
$ws.Range("J71").Value = 3431.3809
$ws.Range("L71").Value = 30882.4281
$ws.Range("N71").Value = -38994.4281

$ws.Range("H107").Value = 3404.5
$ws.Range("I107").Value = 2531
$ws.Range("J107").Value = 3841.25
$ws.Range("K107").Value = 7593
$ws.Range("L107").Value = 11523.75
$ws.Range("M107").Value = -5673
$ws.Range("N107").Value = -15363.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4092
$ws.Range("I80").Value = 2851.6667
$ws.Range("J80").Value = 5952.5
$ws.Range("K80").Value = 2851.6667
$ws.Range("L80").Value = 5952.5
$ws.Range("M80").Value = -1853.6667
$ws.Range("N80").Value = -7948.5

$ws.Range("H83").Value = 4092
$ws.Range("I83").Value = 2851.6667
$ws.Range("J83").Value = 5952.5
$ws.Range("K83").Value = 14258.3335
$ws.Range("L83").Value = 29762.5
$ws.Range("M83").Value = -9266.333500000001
$ws.Range("N83").Value = -39746.5

$ws.Range("H132").Value = 5993.2607
$ws.Range("I132").Value = 5991.222
$ws.Range("J132").Value = 6000.6
$ws.Range("K132").Value = 17973.666
$ws.Range("L132").Value = 18001.8
$ws.Range("M132").Value = -15443.666
$ws.Range("N132").Value = -23061.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5831.905
$ws.Range("I22").Value = 1472.75
$ws.Range("J22").Value = 6857.5884
$ws.Range("K22").Value = 1472.75
$ws.Range("L22").Value = 6857.5884
$ws.Range("M22").Value = -1177.75
$ws.Range("N22").Value = -7447.5884

$ws.Range("H27").Value = 5831.905
$ws.Range("I27").Value = 1472.75
$ws.Range("J27").Value = 6857.5884
$ws.Range("K27").Value = 1472.75
$ws.Range("L27").Value = 6857.5884
$ws.Range("M27").Value = -1365.75
$ws.Range("N27").Value = -7071.5884

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 982.94116
$ws.Range("I132").Value = 919.375
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2758.125
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -228.125
$ws.Range("N132").Value = -11060
